$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '24.957.71'
$ws.Range('E2').Value = '  -3.91%  '
$ws.Range('D3').Value = '1.638.56'
$ws.Range('E3').Value = '  -6.13%  '
$ws.Range('D4').Value = '0.9973'
$ws.Range('E4').Value = '  -0.27%  '
$ws.Range('D5').Value = '233.69'
$ws.Range('E5').Value = '  -6.06%  '
$ws.Range('E6').Value = '  -0.11%  '
$ws.Range('D7').Value = '0.4761'
$ws.Range('E7').Value = '  -6.24%  '
$ws.Range('D8').Value = '39.15'
$ws.Range('E8').Value = '  -4.14%  '
$ws.Range('E9').Value = '  -6.01%  '
$ws.Range('D10').Value = '0.06087'
$ws.Range('E10').Value = '  -1.73%  '
$ws.Range('D11').Value = '0.07016'
$ws.Range('E11').Value = '  -3.22%  '
$ws.Range('D12').Value = '1.640.14'
$ws.Range('E12').Value = '  -6.01%  '
$ws.Range('D13').Value = '14.65'
$ws.Range('E13').Value = '  -3.60%  '
$ws.Range('D14').Value = '0.6035'
$ws.Range('E14').Value = '  -7.76%  '
$ws.Range('D15').Value = '4.343'
$ws.Range('E15').Value = '  -7.13%  '
$ws.Range('D16').Value = '''73.50'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -5.50%  '
$ws.Range('E17').Value = '  -0.12%  '
$ws.Range('D18').Value = '0.9979'
$ws.Range('E18').Value = '  -0.19%  '
$ws.Range('D19').Value = '24.964.09'
$ws.Range('D20').Value = '''0.000006588'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.99%  '
$ws.Range('E21').Value = '  -5.80%  '
$ws.Range('D22').Value = '1.850.52'
$ws.Range('E22').Value = '  -6.07%  '
$ws.Range('D23').Value = '4.377'
$ws.Range('E23').Value = '  -1.22%  '
$ws.Range('D24').Value = '8.586'
$ws.Range('E24').Value = '  -1.87%  '
$ws.Range('D25').Value = '5.269'
$ws.Range('E25').Value = '  -2.23%  '
$ws.Range('D26').Value = '133.52'
$ws.Range('E26').Value = '  -2.30%  '
$ws.Range('D27').Value = '14.85'
$ws.Range('E27').Value = '  -2.73%  '
$ws.Range('D28').Value = '1.385'
$ws.Range('E28').Value = '  -8.81%  '
$ws.Range('D29').Value = '103.69'
$ws.Range('E29').Value = '  -1.98%  '
$ws.Range('E30').Value = '  -8.18%  '
$ws.Range('D31').Value = '3.949'
$ws.Range('E31').Value = '  +2.33%  '
$ws.Range('D32').Value = '0.07721'
$ws.Range('E32').Value = '  -5.87%  '
$ws.Range('D33').Value = '3.548'
$ws.Range('E33').Value = '  -2.77%  '
$ws.Range('D34').Value = '0.9978'
$ws.Range('E34').Value = '  -0.13%  '
$ws.Range('D35').Value = '0.04301'
$ws.Range('E35').Value = '  -8.02%  '
$ws.Range('D36').Value = '2.588'
$ws.Range('D37').Value = '0.9259'
$ws.Range('E37').Value = '  -7.22%  '
$ws.Range('D38').Value = '0.5838'
$ws.Range('E38').Value = '  -5.42%  '
$ws.Range('D39').Value = '2.539'
$ws.Range('E39').Value = '  -7.88%  '
$ws.Range('D40').Value = '0.01534'
$ws.Range('E40').Value = '  -5.12%  '
$ws.Range('D41').Value = '0.9977'
$ws.Range('E41').Value = '  -0.21%  '
$ws.Range('D42').Value = '0.8228'
$ws.Range('E42').Value = '  +7.16%  '
$ws.Range('D43').Value = '97.83'
$ws.Range('E43').Value = '  -2.98%  '
$ws.Range('D44').Value = '1.765'
$ws.Range('E44').Value = '  -8.73%  '
$ws.Range('D45').Value = '0.3699'
$ws.Range('E45').Value = '  -5.90%  '
$ws.Range('D46').Value = '4.698'
$ws.Range('E46').Value = '  -6.12%  '
$ws.Range('B47').Value = 'Aptos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D47').Value = '6.085'
$ws.Range('E47').Value = '  -4.21%  '
$ws.Range('B48').Value = 'Algorand'
$ws.Range('C48').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D48').Value = '0.1089'
$ws.Range('E48').Value = '  -5.62%  '
$ws.Range('E49').Value = '  -2.81%  '
$ws.Range('E50').Value = '  -4.12%  '
$ws.Range('D51').Value = '0.9984'
$ws.Range('E51').Value = '  -0.40%  '
